# The design doc header/banner area was shrunk, so the body content shapes on
# slide 1 shift upward (Top decreases, Left unchanged). On slide 3 a caption
# textbox and a picture near the bottom-left are also repositioned slightly.
$p = $ppt.ActivePresentation

# --- Slide 1 ---
$s = $p.Slides.Item(1)
$moves = @(
    @(1, 84.65937, 166.945199),
    @(2, 243.191582, 177.626381),
    @(3, 99.751973, 177.626381),
    @(4, 348.923081, 45.099135),
    @(5, 357.758194, 69.23087),
    @(6, 470.835129, 69.23087),
    @(7, 528.593307, 69.23087),
    @(8, 600.199213, 69.23087),
    @(9, 576.562531, 71.701811),
    @(10, 357.758194, 113.077008),
    @(11, 432.494339, 113.077008),
    @(12, 539.027741, 111.165906),
    @(13, 511.667717, 112.844174),
    @(14, 407.934016, 94.154019),
    @(15, 432.494339, 94.154019),
    @(16, 462.922928, 94.154019),
    @(17, 375.725118, 200.571654),
    @(18, 367.912126, 138.000161),
    @(19, 422.175919, 94.154019),
    @(20, 492.050315, 94.153937),
    @(21, 364.523071, 269.107956),
    @(22, 351.667008, 42.896693),
    @(23, 697.535512, 46.849608),
    @(24, 704.634675, 42.896693),
    @(25, 718.211181, 69.477087),
    @(26, 776.82071, 69.477087),
    @(27, 833.331654, 68.947247),
    @(28, 888.693228, 69.477087),
    @(29, 712.742283, 108.716694),
    @(30, 796.096955, 108.716694),
    @(31, 865.605039, 108.716694),
    @(32, 435.81622, 282.593622),
    @(33, 365.475827, 325.292441),
    @(34, 600.182047, 111.90811),
    @(35, 535.725512, 138.000161),
    @(36, 519.494355, 136.089058),
    @(37, 459.326615, 138.000161),
    @(38, 255.478035, 233.700005),
    @(39, 253.460717, 270.476457),
    @(40, 243.191582, 176.044094),
    @(41, 590.242524, 257.934252),
    @(42, 243.191582, 389.343465),
    @(43, 408.456223, 356.373551),
    @(44, 243.420952, 389.956223),
    @(45, 591.962441, 395.208976),
    @(46, 590.242524, 389.696063),
    @(47, 773.041596, 246.091103),
    @(48, 776.82071, 380.646299),
    @(49, 244.060952, 433.040157),
    @(50, 244.060952, 433.137406),
    @(51, 548.209687, 131.848661),
    @(52, 565.500641, 133.412209),
    @(53, 99.751973, 247.924889),
    @(54, 763.212128, 221.129134),
    @(55, 761.807038, 355.72316),
    @(56, 156.718976, 207.692362),
    @(57, 5.695749, 45.099135),
    @(58, 213.686063, 262.32189),
    @(59, 205.140008, 284.985039),
    @(60, 183.121895, 281.946772),
    @(61, 242.657953, 469.614579),
    @(62, 242.58189, 471.379055),
    @(63, 788.593071, 436.354966),
    @(64, -15.198977, 148.00685),
    @(65, 788.593071, 483.556229),
    @(66, 587.807244, 450.286772),
    @(67, 587.807244, 487.94378)
)
foreach ($m in $moves) {
    $sh = $s.Shapes.Item([int]$m[0])
    $sh.Left = $m[1]
    $sh.Top = $m[2]
}

# --- Slide 3 ---
$s = $p.Slides.Item(3)
$moves = @(
    @(5, 175.464806, 415.15828),
    @(21, 137.633386, 414.346695)
)
foreach ($m in $moves) {
    $sh = $s.Shapes.Item([int]$m[0])
    $sh.Left = $m[1]
    $sh.Top = $m[2]
}

